# A number of MAG entries were dropped from the "o__RF32_pred-t-p" sheet
# (leftover rows from an out-of-date / renamed input set). Remove those
# specific data rows so the remaining rows shift up and the sheet's used
# range shrinks from A1:F160 to A1:F134, matching the refreshed output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original (pre-edit) row numbers of the entries to remove, listed from
# highest to lowest so deleting one doesn't change the row numbers of the
# others still waiting to be deleted.
$rowsToDelete = @(148, 144, 142, 141, 127, 110, 104, 97, 96, 90, 87, 85, 84, 83, 76, 75, 71, 67, 65, 52, 30, 26, 24, 12, 8, 3)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
